$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 529.6667
$ws.Range("I18").Value = 470.875
$ws.Range("K18").Value = 470.875
$ws.Range("M18").Value = -186.875
$ws.Range("H40").Value = 1601.3334
$ws.Range("I40").Value = 1529
$ws.Range("J40").Value = 1637.5
$ws.Range("K40").Value = 1529
$ws.Range("L40").Value = 1637.5
$ws.Range("M40").Value = -1354
$ws.Range("N40").Value = -1987.5
$ws.Range("H51").Value = 4000
$ws.Range("J51").Value = 4000
$ws.Range("L51").Value = 4000
$ws.Range("N51").Value = -4968
$ws.Range("H62").Value = 996.5
$ws.Range("I62").Value = 996.5
$ws.Range("K62").Value = 996.5
$ws.Range("M62").Value = -372.5
$ws.Range("H65").Value = 996.5
$ws.Range("I65").Value = 996.5
$ws.Range("K65").Value = 4982.5
$ws.Range("M65").Value = -1862.5
$ws.Range("H97").Value = 1150
$ws.Range("J97").Value = 1150
$ws.Range("L97").Value = 3450
$ws.Range("N97").Value = -4442
$ws.Range("H113").Value = 5841.857
$ws.Range("I113").Value = 5982.1665
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 5982.1665
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2728.1665
$ws.Range("N113").Value = -11508
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 1065.5
$ws.Range("I132").Value = 1123.25
$ws.Range("J132").Value = 950
$ws.Range("K132").Value = 3369.75
$ws.Range("L132").Value = 2850
$ws.Range("M132").Value = -839.75
$ws.Range("N132").Value = -7910
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 30000
$ws.Range("K44").Value = 30000
$ws.Range("M44").Value = -29512
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 257.625
$ws.Range("I5").Value = 232.2
$ws.Range("K5").Value = 232.2
$ws.Range("M5").Value = -119.2
$ws.Range("H20").Value = 1105.8334
$ws.Range("I20").Value = 850.3333
$ws.Range("J20").Value = 1361.3334
$ws.Range("K20").Value = 850.3333
$ws.Range("L20").Value = 1361.3334
$ws.Range("M20").Value = -603.3333
$ws.Range("N20").Value = -1855.3334
$ws.Range("H36").Value = 4499.4
$ws.Range("I36").Value = 4499.4
$ws.Range("K36").Value = 4499.4
$ws.Range("M36").Value = -3965.4
$ws.Range("H86").Value = 2100
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -77
$ws.Range("N86").Value = -5246
$ws.Range("H88").Value = 18351.908
$ws.Range("J88").Value = 18351.908
$ws.Range("L88").Value = 18351.908
$ws.Range("N88").Value = -19163.908
$ws.Range("H89").Value = 2100
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -384
$ws.Range("N89").Value = -26232
$ws.Range("H91").Value = 18351.908
$ws.Range("J91").Value = 18351.908
$ws.Range("L91").Value = 18351.908
$ws.Range("N91").Value = -21159.908
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1137.4286
$ws.Range("I2").Value = 891.9
$ws.Range("J2").Value = 1751.25
$ws.Range("K2").Value = 891.9
$ws.Range("L2").Value = 1751.25
$ws.Range("M2").Value = -778.9
$ws.Range("N2").Value = -1977.25
$ws.Range("H16").Value = 2219.4614
$ws.Range("I16").Value = 1425.4
$ws.Range("J16").Value = 4866.3335
$ws.Range("K16").Value = 1425.4
$ws.Range("L16").Value = 4866.3335
$ws.Range("M16").Value = -1138.4
$ws.Range("N16").Value = -5440.3335
$ws.Range("H36").Value = 25600
$ws.Range("I36").Value = 14500
$ws.Range("K36").Value = 14500
$ws.Range("M36").Value = -14112
$ws.Range("H40").Value = 25600
$ws.Range("I40").Value = 14500
$ws.Range("K40").Value = 14500
$ws.Range("M40").Value = -14340
$ws.Range("H55").Value = 43333.332
$ws.Range("I55").Value = 45000
$ws.Range("K55").Value = 45000
$ws.Range("M55").Value = -44685
$ws.Range("H99").Value = 1432142.1
$ws.Range("I99").Value = 2002999.4
$ws.Range("K99").Value = 2002999.4
$ws.Range("M99").Value = -2001501.4
$ws.Range("H113").Value = 2219.4614
$ws.Range("I113").Value = 1425.4
$ws.Range("J113").Value = 4866.3335
$ws.Range("K113").Value = 1425.4
$ws.Range("L113").Value = 4866.3335
$ws.Range("M113").Value = 744.5999999999999
$ws.Range("N113").Value = -9206.333500000001
$ws.Range("H122").Value = 3974.125
$ws.Range("I122").Value = 915.5
$ws.Range("J122").Value = 7032.75
$ws.Range("K122").Value = 2746.5
$ws.Range("L122").Value = 21098.25
$ws.Range("M122").Value = -296.5
$ws.Range("N122").Value = -25998.25
$ws.Range("H126").Value = 1432142.1
$ws.Range("I126").Value = 2002999.4
$ws.Range("K126").Value = 6008998.199999999
$ws.Range("M126").Value = -6006528.199999999
$ws.Range("H132").Value = 2975.5
$ws.Range("I132").Value = 2975.5
$ws.Range("K132").Value = 8926.5
$ws.Range("M132").Value = -6396.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H34").Value = 14227
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H44").Value = 529.4
$ws.Range("J44").Value = 851.6
$ws.Range("L44").Value = 2554.8
$ws.Range("N44").Value = -3350.8
$ws.Range("H50").Value = 497.25
$ws.Range("I50").Value = 497.25
$ws.Range("K50").Value = 1491.75
$ws.Range("M50").Value = -1010.75
$ws.Range("H53").Value = 497.25
$ws.Range("I53").Value = 497.25
$ws.Range("K53").Value = 1491.75
$ws.Range("M53").Value = -1010.75
$ws.Range("H55").Value = 2505.303
$ws.Range("I55").Value = 1218.1818
$ws.Range("J55").Value = 3148.8635
$ws.Range("K55").Value = 3654.5454
$ws.Range("L55").Value = 9446.5905
$ws.Range("M55").Value = -3477.5454
$ws.Range("N55").Value = -9800.5905
$ws.Range("H68").Value = 732.875
$ws.Range("I68").Value = 697
$ws.Range("K68").Value = 2091
$ws.Range("M68").Value = -1280
$ws.Range("H71").Value = 732.875
$ws.Range("I71").Value = 697
$ws.Range("K71").Value = 6273
$ws.Range("M71").Value = -2217
$ws.Range("H97").Value = 518.75
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 275
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 825
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -1817
$ws.Range("H140").Value = 638.6667
$ws.Range("I140").Value = 638.6667
$ws.Range("K140").Value = 1916.0001
$ws.Range("M140").Value = 3263.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 4940
$ws.Range("J98").Value = 4940
$ws.Range("L98").Value = 4940
$ws.Range("N98").Value = -10930
$ws.Range("H113").Value = 5624.875
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1064
$ws.Range("J22").Value = 1196
$ws.Range("L22").Value = 1196
$ws.Range("N22").Value = -1786
$ws.Range("H27").Value = 1064
$ws.Range("J27").Value = 1196
$ws.Range("L27").Value = 1196
$ws.Range("N27").Value = -1410
$ws.Range("H40").Value = 390955.3
$ws.Range("J40").Value = 842317.8
$ws.Range("L40").Value = 842317.8
$ws.Range("N40").Value = -842589.8
$ws.Range("H46").Value = 5500
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812
$ws.Range("H55").Value = 660.1667
$ws.Range("I55").Value = 592.2
$ws.Range("K55").Value = 592.2
$ws.Range("M55").Value = -419.2
$ws.Range("H82").Value = 1542.6364
$ws.Range("I82").Value = 1480.6666
$ws.Range("J82").Value = 1617
$ws.Range("K82").Value = 1480.6666
$ws.Range("L82").Value = 1617
$ws.Range("M82").Value = -1119.6666
$ws.Range("N82").Value = -2339
$ws.Range("H85").Value = 1542.6364
$ws.Range("I85").Value = 1480.6666
$ws.Range("J85").Value = 1617
$ws.Range("K85").Value = 1480.6666
$ws.Range("L85").Value = 1617
$ws.Range("M85").Value = -232.6666
$ws.Range("N85").Value = -4113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1967.6
$ws.Range("I96").Value = 2016
$ws.Range("J96").Value = 1895
$ws.Range("K96").Value = 2016
$ws.Range("L96").Value = 1895
$ws.Range("M96").Value = -643
$ws.Range("N96").Value = -4641
$ws.Range("H122").Value = 1433.3214
$ws.Range("I122").Value = 1428.1923
$ws.Range("K122").Value = 4284.5769
$ws.Range("M122").Value = -1834.5769
$ws.Range("H141").Value = 40599.6
$ws.Range("I141").Value = 41999.5
$ws.Range("J141").Value = 35000
$ws.Range("K141").Value = 41999.5
$ws.Range("L141").Value = 35000
$ws.Range("M141").Value = -36819.5
$ws.Range("N141").Value = -45360
